$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.053462380191461
$ws.Range("D2").Value = 1.054942577598969
$ws.Range("E2").Value = 0.992614727750844
$ws.Range("F2").Value = 1.064246972990651
$ws.Range("I2").Value = 1.041217604814249
$ws.Range("J2").Value = 1.058479597006812
$ws.Range("K2").Value = 1.057684398785629
$ws.Range("L2").Value = 0.9955398523335997
$ws.Range("M2").Value = 1.0669634261984
$ws.Range("N2").Value = 1.059982759443869
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.054987353350884
$ws.Range("D3").Value = 1.056084274880118
$ws.Range("E3").Value = 0.9936372048519299
$ws.Range("F3").Value = 1.065537212671073
$ws.Range("I3").Value = 1.041552605815277
$ws.Range("J3").Value = 1.05965267101323
$ws.Range("K3").Value = 1.058638600962229
$ws.Range("L3").Value = 0.9963617723202687
$ws.Range("M3").Value = 1.068067643831
$ws.Range("N3").Value = 1.061157499349929
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.055972533699497
$ws.Range("D4").Value = 1.056821532110635
$ws.Range("E4").Value = 0.9942998659930998
$ws.Range("F4").Value = 1.066370825462873
$ws.Range("I4").Value = 1.041767382086344
$ws.Range("J4").Value = 1.060409754689269
$ws.Range("K4").Value = 1.059253957113985
$ws.Range("L4").Value = 0.9968940712668347
$ws.Range("M4").Value = 1.068780336261158
$ws.Range("N4").Value = 1.061915658171627
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.056386333286141
$ws.Range("D5").Value = 1.05713112142952
$ws.Range("E5").Value = 0.994578699834602
$ws.Range("F5").Value = 1.066720980364802
$ws.Range("I5").Value = 1.04185719905393
$ws.Range("J5").Value = 1.06072756595737
$ws.Range("K5").Value = 1.059512160284309
$ws.Range("L5").Value = 0.9971179600053012
$ws.Range("M5").Value = 1.069079523772292
$ws.Range("N5").Value = 1.062233920768182
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.056455790574099
$ws.Range("D6").Value = 1.057183082241479
$ws.Range("E6").Value = 0.994625531979634
$ws.Range("F6").Value = 1.066779755716765
$ws.Range("I6").Value = 1.041872251918872
$ws.Range("J6").Value = 1.060780900667083
$ws.Range("K6").Value = 1.059555484980748
$ws.Range("L6").Value = 0.9971555583673455
$ws.Range("M6").Value = 1.06912973369525
$ws.Range("N6").Value = 1.0622873312193
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.055978064354887
$ws.Range("D7").Value = 1.056825670241806
$ws.Range("E7").Value = 0.994303590798249
$ws.Range("F7").Value = 1.066375505405908
$ws.Range("I7").Value = 1.04176858408911
$ws.Range("J7").Value = 1.060414003125701
$ws.Range("K7").Value = 1.059257409167184
$ws.Range("L7").Value = 0.9968970624462089
$ws.Range("M7").Value = 1.068784335696448
$ws.Range("N7").Value = 1.061919912641326
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.053978083741795
$ws.Range("D8").Value = 1.055328732471011
$ws.Range("E8").Value = 0.9929600610674297
$ws.Range("F8").Value = 1.064683279376441
$ws.Range("I8").Value = 1.041331233504913
$ws.Range("J8").Value = 1.058876455009911
$ws.Range("K8").Value = 1.05800730857421
$ws.Range("L8").Value = 0.9958175282591056
$ws.Range("M8").Value = 1.067336979543547
$ws.Range("N8").Value = 1.060380181030852
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.050441414425814
$ws.Range("D9").Value = 1.052679252262559
$ws.Range("E9").Value = 0.9906006454969559
$ws.Range("F9").Value = 1.061691479296008
$ws.Range("I9").Value = 1.040545227363666
$ws.Range("J9").Value = 1.056151722566762
$ws.Range("K9").Value = 1.055788363230371
$ws.Range("L9").Value = 0.9939188001724441
$ws.Range("M9").Value = 1.064772464533508
$ws.Range("N9").Value = 1.057651579155103
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.048074770127937
$ws.Range("D10").Value = 1.05090476926654
$ws.Range("E10").Value = 0.989033133672735
$ws.Range("F10").Value = 1.05968996008534
$ws.Range("I10").Value = 1.040010791774625
$ws.Range("J10").Value = 1.054324529542675
$ws.Range("K10").Value = 1.054297946645234
$ws.Range("L10").Value = 0.9926553831429383
$ws.Range("M10").Value = 1.063052999019476
$ws.Range("N10").Value = 1.05582179130734
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.047047773506373
$ws.Range("D11").Value = 1.050134393318796
$ws.Range("E11").Value = 0.988355674866747
$ws.Range("F11").Value = 1.058821546283757
$ws.Range("I11").Value = 1.039776874190313
$ws.Range("J11").Value = 1.053530713451988
$ws.Range("K11").Value = 1.05364987724555
$ws.Range("L11").Value = 0.9921088820399291
$ws.Range("M11").Value = 1.062306058465713
$ws.Range("N11").Value = 1.055026847906751
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.046665957459032
$ws.Range("D12").Value = 1.049847932892689
$ws.Range("E12").Value = 0.9881042295826724
$ws.Range("F12").Value = 1.058498710171599
$ws.Range("I12").Value = 1.039689608289501
$ws.Range("J12").Value = 1.053235452577807
$ws.Range("K12").Value = 1.053408742684251
$ws.Range("L12").Value = 0.9919059725120875
$ws.Range("M12").Value = 1.062028244582434
$ws.Range("N12").Value = 1.054731167728263
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.046747874009443
$ws.Range("D13").Value = 1.049909393684449
$ws.Range("E13").Value = 0.9881581567098651
$ws.Range("F13").Value = 1.058567971875108
$ws.Range("I13").Value = 1.039708344287591
$ws.Range("J13").Value = 1.053298805355563
$ws.Range("K13").Value = 1.053460485635307
$ws.Range("L13").Value = 0.9919494934313052
$ws.Range("M13").Value = 1.062087853307331
$ws.Range("N13").Value = 1.05479461047423
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.047016219541741
$ws.Range("D14").Value = 1.050110720739087
$ws.Range("E14").Value = 0.9883348863814464
$ws.Range("F14").Value = 1.05879486606769
$ws.Range("I14").Value = 1.039769668498756
$ws.Range("J14").Value = 1.053506315343256
$ws.Range("K14").Value = 1.05362995344027
$ws.Range("L14").Value = 0.9920921077337197
$ws.Range("M14").Value = 1.062283101804452
$ws.Range("N14").Value = 1.055002415149907
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.047181510237455
$ws.Range("D15").Value = 1.050234723921087
$ws.Range("E15").Value = 0.9884438009545853
$ws.Range("F15").Value = 1.058934627279097
$ws.Range("I15").Value = 1.039807402143393
$ws.Range("J15").Value = 1.053634115579468
$ws.Range("K15").Value = 1.053734313285605
$ws.Range("L15").Value = 0.9921799884222134
$ws.Range("M15").Value = 1.062403352041724
$ws.Range("N15").Value = 1.055130396877113
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.048142880832627
$ws.Range("D16").Value = 1.05095585365193
$ws.Range("E16").Value = 0.9890781214508737
$ws.Range("F16").Value = 1.059747556558108
$ws.Range("I16").Value = 1.040026263166887
$ws.Range("J16").Value = 1.054377156490545
$ws.Range("K16").Value = 1.054340899311591
$ws.Range("L16").Value = 0.9926916645766087
$ws.Range("M16").Value = 1.063102519894369
$ws.Range("N16").Value = 1.055874492991512
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.04874532099126
$ws.Range("D17").Value = 1.051407656059647
$ws.Range("E17").Value = 0.9894763578477731
$ws.Range("F17").Value = 1.06025701429539
$ws.Range("I17").Value = 1.040162876987103
$ws.Range("J17").Value = 1.054842537388034
$ws.Range("K17").Value = 1.054720665361234
$ws.Range("L17").Value = 0.9930127773692701
$ws.Range("M17").Value = 1.063540442312121
$ws.Range("N17").Value = 1.056340534783255
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.049096500213
$ws.Range("D18").Value = 1.051670990937517
$ws.Range("E18").Value = 0.9897087662937551
$ws.Range("F18").Value = 1.060554004679886
$ws.Range("I18").Value = 1.040242320132978
$ws.Range("J18").Value = 1.055113732977004
$ws.Range("K18").Value = 1.054941915464661
$ws.Range("L18").Value = 0.9932001317071766
$ws.Range("M18").Value = 1.063795643847099
$ws.Range("N18").Value = 1.056612115501068
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.049216207256287
$ws.Range("D19").Value = 1.051760748631983
$ws.Range("E19").Value = 0.9897880325774039
$ws.Range("F19").Value = 1.060655242484544
$ws.Range("I19").Value = 1.040269367321322
$ws.Range("J19").Value = 1.055206160931149
$ws.Range("K19").Value = 1.055017311932426
$ws.Range("L19").Value = 0.993264023964098
$ws.Range("M19").Value = 1.063882621947137
$ws.Range("N19").Value = 1.05670467471351
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.048680707061333
$ws.Range("D20").Value = 1.051359202028989
$ws.Range("E20").Value = 0.9894336180355766
$ws.Range("F20").Value = 1.060202371695835
$ws.Range("I20").Value = 1.04014824460185
$ws.Range("J20").Value = 1.054792632665351
$ws.Range("K20").Value = 1.054679947096849
$ws.Range("L20").Value = 0.9929783193490043
$ws.Range("M20").Value = 1.063493481340689
$ws.Range("N20").Value = 1.056290559190141
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.046937208038494
$ws.Range("D21").Value = 1.050051443504376
$ws.Range("E21").Value = 0.9882828385668255
$ws.Range("F21").Value = 1.058728058851059
$ws.Range("I21").Value = 1.039751620514566
$ws.Range("J21").Value = 1.053445219985567
$ws.Range("K21").Value = 1.053580060835057
$ws.Range("L21").Value = 0.9920501090198107
$ws.Range("M21").Value = 1.062225616167183
$ws.Range("N21").Value = 1.054941233029802
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.045839007196183
$ws.Range("D22").Value = 1.049227415902662
$ws.Range("E22").Value = 0.9875604150241496
$ws.Range("F22").Value = 1.05779954266907
$ws.Range("I22").Value = 1.039500056159782
$ws.Range("J22").Value = 1.052595717352305
$ws.Range("K22").Value = 1.052886127890563
$ws.Range("L22").Value = 0.991467000034148
$ws.Range("M22").Value = 1.061426333089228
$ws.Range("N22").Value = 1.054090524005362
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.046421376469075
$ws.Range("D23").Value = 1.049664420033109
$ws.Range("E23").Value = 0.9879432794636459
$ws.Range("F23").Value = 1.058291916528769
$ws.Range("I23").Value = 1.039633623670204
$ws.Range("J23").Value = 1.053046278119303
$ws.Range("K23").Value = 1.053254223411854
$ws.Range("L23").Value = 0.9917760702887607
$ws.Range("M23").Value = 1.061850251779872
$ws.Range("N23").Value = 1.054541724620329
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.048709903976234
$ws.Range("D24").Value = 1.051381096918924
$ws.Range("E24").Value = 0.9894529299347241
$ws.Range("F24").Value = 1.060227062852959
$ws.Range("I24").Value = 1.040154857092862
$ws.Range("J24").Value = 1.054815183244038
$ws.Range("K24").Value = 1.054698346734531
$ws.Range("L24").Value = 0.9929938892766438
$ws.Range("M24").Value = 1.063514701696906
$ws.Range("N24").Value = 1.056313141793237
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.051357254145445
$ws.Range("D25").Value = 1.053365622507484
$ws.Range("E25").Value = 0.9912096547607046
$ws.Range("F25").Value = 1.062466137381495
$ws.Range("I25").Value = 1.040750258641535
$ws.Range("J25").Value = 1.056857990937748
$ws.Range("K25").Value = 1.056363952347499
$ws.Range("L25").Value = 0.9944092447426411
$ws.Range("M25").Value = 1.065437155117135
$ws.Range("N25").Value = 1.058358850508186
